$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Looking for the Tomcat documentation on Internet, I found that the",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Looking for the Tomcat documentation on the Internet, I found that the",
    2)
